$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.825.09'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '1.861.91'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5074'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3651'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07150'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8886'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07485'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').Value = '1.864.32'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.212'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008474'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.0000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '26.852.56'
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.987'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').Value = '2.107.71'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.342'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.772'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.30%  '
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.083'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.57'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.670'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.693'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09116'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.36%  '
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7463'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.942'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.147'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.209'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.493'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5539'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01974'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.554'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.556'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1483'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4735'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.0000'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.02'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.547'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.73%  '
